$d = $word.ActiveDocument

# 1. Change heading "Built-In Allocator" -> "Customizable Allocator"
$d.Content.Find.Execute("Built-In Allocator", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Customizable Allocator", 2)

# 2. Replace the paragraph body text about the block allocator
$d.Content.Find.Execute("Jinx utilizes its own block allocator designed to prioritize efficiency for small, frequent allocations, as is typical of scripting requirements.  Additionally, it makes use of thread-local storage pools to ensure minimal contention between scripts executing independently on different threads.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Jinx allows the user to supply a custom allocator, potentially enabling better performance than with the default system allocator.", 2)
